$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'30.488.39"
$ws.Range('E2').Value = '  -1.04%  '
$ws.Range('D3').Value = "'2.092.35"
$ws.Range('E3').Value = '  -1.19%  '
$ws.Range('D4').Value = "'1.003"
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = "'330.12"
$ws.Range('E5').Value = '  -1.52%  '
$ws.Range('D6').Value = "'1.002"
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('D7').Value = "'0.5212"
$ws.Range('E7').Value = '  -4.15%  '
$ws.Range('D8').Value = "'0.4418"
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').Value = "'53.51"
$ws.Range('E9').Value = '  +14.68%  '
$ws.Range('D10').Value = "'0.08925"
$ws.Range('E10').Value = '  -1.61%  '
$ws.Range('D11').Value = "'1.152"
$ws.Range('E11').Value = '  -3.27%  '
$ws.Range('D12').Value = "'24.22"
$ws.Range('E12').Value = '  -4.50%  '
$ws.Range('D13').Value = "'2.102.57"
$ws.Range('E13').Value = '  -0.83%  '
$ws.Range('D14').Value = "'6.684"
$ws.Range('E14').Value = '  -1.72%  '
$ws.Range('D15').Value = "'7.676"
$ws.Range('E15').Value = '  -2.42%  '
$ws.Range('D16').Value = "'95.91"
$ws.Range('E16').Value = '  -2.63%  '
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('D18').Value = "'0.00001123"
$ws.Range('E18').Value = '  -1.81%  '
$ws.Range('D19').Value = "'0.06608"
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('E20').Value = '  -0.81%  '
$ws.Range('D21').Value = "'1.001"
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').Value = "'6.256"
$ws.Range('E22').Value = '  -2.60%  '
$ws.Range('D23').Value = "'30.526.54"
$ws.Range('E23').Value = '  -1.27%  '
$ws.Range('D24').Value = "'12.30"
$ws.Range('E24').Value = '  +0.86%  '
$ws.Range('D25').Value = "'2.318"
$ws.Range('E25').Value = '  +1.93%  '
$ws.Range('D26').Value = "'2.345.81"
$ws.Range('E26').Value = '  -0.94%  '
$ws.Range('E27').Value = '  -3.18%  '
$ws.Range('E28').Value = '  -0.87%  '
$ws.Range('D29').Value = "'163.58"
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('D30').Value = "'131.48"
$ws.Range('E30').Value = '  -1.94%  '
$ws.Range('D31').Value = "'1.189"
$ws.Range('E31').Value = '  +1.18%  '
$ws.Range('E32').Value = '  -1.87%  '
$ws.Range('D33').Value = "'1.655"
$ws.Range('E33').Value = '  +6.43%  '
$ws.Range('D34').Value = "'6.150"
$ws.Range('E34').Value = '  -2.45%  '
$ws.Range('E35').Value = '  -2.89%  '
$ws.Range('D36').Value = "'10.01"
$ws.Range('E36').Value = '  +3.78%  '
$ws.Range('D37').Value = "'0.02565"
$ws.Range('E37').Value = '  -1.65%  '
$ws.Range('D38').Value = "'0.06827"
$ws.Range('E38').Value = '  +0.93%  '
$ws.Range('D39').Value = "'5.474"
$ws.Range('E39').Value = '  -2.14%  '
$ws.Range('D40').Value = "'12.58"
$ws.Range('E40').Value = '  -1.52%  '
$ws.Range('E41').Value = '  -1.60%  '
$ws.Range('D42').Value = "'0.6876"
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('D43').Value = "'1.249"
$ws.Range('E43').Value = '  -1.11%  '
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('D45').Value = "'13.99"
$ws.Range('E45').Value = '  -1.10%  '
$ws.Range('D46').Value = "'0.6325"
$ws.Range('E46').Value = '  -2.35%  '
$ws.Range('D47').Value = "'2.195"
$ws.Range('E47').Value = '  -3.25%  '
$ws.Range('D48').Value = "'3.628"
$ws.Range('E48').Value = '  -1.54%  '
$ws.Range('E49').Value = '  +6.08%  '
$ws.Range('D50').Value = "'1.242"
$ws.Range('E50').Value = '  -3.81%  '
$ws.Range('D51').Value = "'81.67"
$ws.Range('E51').Value = '  -2.19%  '
